$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 201 — shifts existing rows 201..281 down to 202..282,
# matching the rest of the table's row layout (A:R).
$ws.Rows.Item(201).Insert()

# Populate the newly-inserted row 201 with the new record.
$ws.Range("A201").Value = 9
$ws.Range("B201").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C201").Value = "Metropolitana"
$ws.Range("D201").Value = 44704
$ws.Range("D201").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E201").Value = 13
$ws.Range("F201").Value = 100112043
$ws.Range("G201").Value = "Pepino ensalada"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 230
$ws.Range("K201").Value = 18000
$ws.Range("L201").Value = 20000
$ws.Range("M201").Value = 18870
$ws.Range("N201").Value = "`$/caja 50 unidades"
$ws.Range("O201").Value = "Región de Arica y Parinacota"
$ws.Range("P201").Value = 377
$ws.Range("Q201").Value = 50
$ws.Range("R201").Value = "Hortaliza"
